# Likitha Seriti WorkFlow Automation.
# Fills in the MarquisDataForm with a new client / vehicle finance deal
# (Hermia Nkosi), replacing the previous placeholder test data, and adds
# a block of new vehicle/finance fields (rows 28-39) including a
# mailto hyperlink on the client-email cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Top "Add New Client Details" block (L/M column, rows 2-16)
# ---------------------------------------------------------------

# First / Last name
$ws.Range("M2").Value = "HERMIA"
$ws.Range("M3").Value = "NKOSI"

# Client Title -> stored as text "MR" (was " MR")
$ws.Range("M4").NumberFormat = "@"
$ws.Range("M4").Value = "MR"

# ID Type
$ws.Range("M5").Value = "RSA ID"

# Finance/Cash Deal (G/H block)
$ws.Range("H6").Value = "FINANCE"

# ID Number -> text (was a plain number)
$ws.Range("M6").NumberFormat = "@"
$ws.Range("M6").Value = "4305045052080"

# DOB -> cleared out, keep date formatting
$ws.Range("M7").ClearContents()

# Citizenship
$ws.Range("M8").Value = "SOUTH AFRICA"

# Mobile Number / Tele code / Tele No -> text (keep leading zeros)
$ws.Range("M10").NumberFormat = "@"
$ws.Range("M10").Value = "0875555555"

$ws.Range("M11").NumberFormat = "@"
$ws.Range("M11").Value = "080"

$ws.Range("M12").NumberFormat = "@"
$ws.Range("M12").Value = "0600777"

# Mob contract type
$ws.Range("M13").Value = "PREPAID"

# ---------------------------------------------------------------
# Employment block (rows 17-27) - only the employment start date
# and the salary figures (now stored as text) change.
# ---------------------------------------------------------------

# Curr emp start date -> new date, keep date formatting
$ws.Range("M22").Value = [DateTime]"1991-05-23"

# Salary Day / Gross Salary
$ws.Range("M23").NumberFormat = "@"
$ws.Range("M23").Value = "25"

$ws.Range("M24").NumberFormat = "@"
$ws.Range("M24").Value = "500000"

# stray formatted (but empty) cell that appears alongside the Net Salary row
$ws.Range("J25").NumberFormat = "@"

$ws.Range("M25").NumberFormat = "@"
$ws.Range("M25").Value = "350000"

# ---------------------------------------------------------------
# New rows 28-39: client email (with hyperlink) + vehicle / finance
# details.
# ---------------------------------------------------------------

$ws.Range("L28").Value = "client email"
$ws.Range("M28").NumberFormat = "@"
$ws.Hyperlinks.Add($ws.Range("M28"), "mailto:tester123@gmail.com", "", "", "tester123@gmail.com")

$ws.Range("L29").Value = "Vehicle Chasis"
$ws.Range("M29").NumberFormat = "@"
$ws.Range("M29").Value = "123456789123356"

$ws.Range("L30").Value = "EngineNumber"
$ws.Range("M30").NumberFormat = "@"
$ws.Range("M30").Value = "23421"

$ws.Range("L31").Value = "Registration Number"
$ws.Range("M31").NumberFormat = "@"
$ws.Range("M31").HorizontalAlignment = -4131
$ws.Range("M31").VerticalAlignment = -4160
$ws.Range("M31").Value = "2332"

$ws.Range("L32").Value = "Vehicle Code"
$ws.Range("M32").NumberFormat = "@"
$ws.Range("M32").Value = "00815170"

$ws.Range("L33").Value = "sellingPrice"
$ws.Range("M33").NumberFormat = "@"
$ws.Range("M33").Value = "20000"

$ws.Range("L34").Value = "firstRegDate"
$ws.Range("M34").NumberFormat = "@"
$ws.Range("M34").Value = "02 Dec 2016"

$ws.Range("L35").Value = "interestRate"
$ws.Range("M35").NumberFormat = "@"
$ws.Range("M35").HorizontalAlignment = -4131
$ws.Range("M35").VerticalAlignment = -4160
$ws.Range("M35").Value = "12.5"

$ws.Range("L36").Value = "vehicle Condition"
$ws.Range("M36").NumberFormat = "@"
$ws.Range("M36").Value = "NEW"

$ws.Range("L37").Value = "Odometer"
$ws.Range("M37").NumberFormat = "@"
$ws.Range("M37").Value = "200"

$ws.Range("L38").Value = "Vehicle Use"
$ws.Range("M38").NumberFormat = "@"
$ws.Range("M38").Value = "PRIVATE"

$ws.Range("L39").Value = "Agreement Type"
$ws.Range("M39").NumberFormat = "@"
$ws.Range("M39").Value = "Instalment Vlaue"

# ---------------------------------------------------------------
# Cosmetic: widen column M a bit to fit the new content, and move
# the selection down to the newly entered row.
# ---------------------------------------------------------------
$ws.Columns("M").ColumnWidth = 20.3

$ws.Range("M39").Select()
